$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab query (B2): remove the trailing Cohort coalesce clause
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Vagina']`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# Update the SamplesTab query (B3): re-entered with no content change
$ws.Range("B3").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) `nWHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Vagina']`n WITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed , `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``"

# Update the FilesTab query (B4): drop the stray blank line before the second WITH DISTINCT
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Vagina']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(f.file_type, '') AS ``File Type``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``File Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed , `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Row heights shrink by one wrapped line for the edited/retyped cells
$ws.Rows(2).RowHeight = 304.5
$ws.Rows(3).RowHeight = 275.5
$ws.Rows(4).RowHeight = 290

# Move the selection/cursor to B2 (matches saved cursor position)
$ws.Range("B2").Select()
